{"js": "// Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific resume bullet\n// paragraphs, matching the target diff exactly.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: the *exact* full paragraph text (used to unambiguously locate\n// the paragraph, since several paragraphs share overlapping prefixes/substrings)\n// plus the ordered list of numeric/metric substrings within it that must be\n// bolded + colored. Order matters only for readability; each token is searched\n// for independently within the paragraph's own range so it cannot bleed into\n// a different paragraph.\nconst EDITS = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    tokens: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    tokens: [\"1,200\"],\n  },\n  {\n    text:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    tokens: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const edit of EDITS) {\n  // Find the (single) paragraph whose full text matches exactly -- plain\n  // substring/body-wide search is unsafe here because some of these target\n  // strings are literal prefixes of other (untouched) paragraphs.\n  let paragraph = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === edit.text) {\n      paragraph = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!paragraph) {\n    continue;\n  }\n\n  for (const token of edit.tokens) {\n    const found = paragraph.search(token, { matchCase: true });\n    found.load(\"text\");\n    await context.sync();\n\n    for (let i = 0; i < found.items.length; i++) {\n      found.items[i].font.bold = true;\n      found.items[i].font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific resume bullet\n# paragraphs, matching the target diff exactly.\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n$pm = [char]0x00B1\n\n# Each entry:\n#   text   = the *exact* full paragraph text (sans trailing paragraph mark),\n#            used to unambiguously locate the paragraph -- several of these\n#            target strings are literal prefixes/substrings of OTHER\n#            (untouched) paragraphs, so plain Find.Execute over the whole\n#            document body is not safe for disambiguation.\n#   tokens = ordered list of numeric/metric substrings within that paragraph\n#            that must be bolded + colored. Each token is located with a\n#            Find.Execute() scoped to a Range clamped to the target\n#            paragraph's own Start/End, so it cannot match text belonging to\n#            a different paragraph.\n$edits = @(\n    @{\n        text = ($bullet + ' Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%')\n        tokens = @('23%', '64%')\n    },\n    @{\n        text = ($bullet + ' Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ' + $pm + '4.2% to ' + $pm + '2.1%')\n        tokens = @('87%', '71%', ($pm + '4.2%'), ($pm + '2.1%'))\n    },\n    @{\n        text = ($bullet + ' Wrote RFP and analyzed bids from 1,200 vendors for research platform development')\n        tokens = @('1,200')\n    },\n    @{\n        text = ($bullet + ' Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+')\n        tokens = @('$400M', '$1B')\n    },\n    @{\n        text = ($bullet + ' Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M')\n        tokens = @('73.5%', '$4.7M')\n    },\n    @{\n        text = ($bullet + ' Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%')\n        tokens = @('87%', '71%')\n    }\n)\n\nforeach ($edit in $edits) {\n    $targetParagraph = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n        if ($t -eq $edit.text) {\n            $targetParagraph = $d.Paragraphs.Item($i)\n            break\n        }\n    }\n    if ($targetParagraph -eq $null) {\n        continue\n    }\n\n    $pStart = $targetParagraph.Range.Start\n    $pEnd = $targetParagraph.Range.End\n\n    foreach ($tok in $edit.tokens) {\n        $sub = $d.Range($pStart, $pEnd)\n        $ok = $sub.Find.Execute($tok)\n        if ($ok) {\n            $sub.Font.Bold = $true\n            $sub.Font.Color = \"#2C3E50\"\n        }\n    }\n}\n"}
